$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C holds the "Förändrad" (changed) date, stored as a serial date
# number. Rows 2-52 all currently hold 45181 (2023-09-12) and need to be
# bumped to 45182 (2023-09-13).
$ws.Range("C2:C52").Value = 45182
